# Add 2022-Q3 data:
#  1. Insert a new row into the "总计" (Total) summary sheet for 2022-Q3
#     and shift the existing rows down.
#  2. Insert a brand-new worksheet named "2022-Q3" right after "总计"
#     (i.e. before the sheet that is currently named "2022-Q2"),
#     populated with the per-fund detail data for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet (first sheet)
# ---------------------------------------------------------------------
# Note: column A here is simply the 0-based row position, so it is
# untouched for rows 2-8 (still 0,1,2,...,6) and a brand-new row 9 is
# appended (value 7). Columns B:D (the quarter label / counts / value)
# get the new 2022-Q3 figures in row 2 and the older quarters shift
# down by one row.
$summary = $wb.Worksheets.Item(1)

# Row 9 is new - give it the same bold/centered/bordered style as the
# rest of column A by copying the format from the cell directly above.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122) # xlPasteFormats
$summary.Range("A9").Value = 7

$quarters = @(
    @("2022-Q3", 9, 0.11),
    @("2022-Q2", 15, 1.63),
    @("2022-Q1", 2, 0.04),
    @("2021-Q4", 4, 0.51),
    @("2021-Q3", 1, 0.03),
    @("2021-Q2", 1, 0.03),
    @("2021-Q1", 4, 1.26),
    @("2020-Q4", 6, 3.24)
)

$r = 2
foreach ($q in $quarters) {
    $summary.Range("B$r").Value = $q[0]
    $summary.Range("C$r").Value = $q[1]
    $summary.Range("D$r").Value = $q[2]
    $r++
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)          # currently "2022-Q2"
$newSheet = $wb.Worksheets.Add($refSheet)   # inserted right before it
$newSheet.Name = "2022-Q3"

# Match page margins used by the other worksheets in this workbook.
$newSheet.PageSetup.LeftMargin = 54     # 0.75 in
$newSheet.PageSetup.RightMargin = 54    # 0.75 in
$newSheet.PageSetup.TopMargin = 72      # 1 in
$newSheet.PageSetup.BottomMargin = 72   # 1 in
$newSheet.PageSetup.HeaderMargin = 36   # 0.5 in
$newSheet.PageSetup.FooterMargin = 36   # 0.5 in

# Header row - copy the bold/centered/bordered style used on every
# other sheet's header row (e.g. from the "总计" sheet) then set text.
$summary.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char](66 + $i)   # B, C, D, E, F, G, H
    $newSheet.Range("$col" + "1").Value = $headers[$i]
}

# Column A (row index) uses the same bold/centered/bordered style too.
$newSheet.Range("B1").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)

$data = @(
    @(0, "519959", "长信多利灵活配置混合A",               "0.95", "89.23", "3.80", "0.0361", 10),
    @(1, "004351", "汇丰晋信珠三角区域发展混合",           "0.42", "93.94", "4.56", "0.0192", 2),
    @(2, "003359", "大成中证360互联网+大数据100指数C",     "1.11", "92.17", "0.98", "0.0109", 9),
    @(3, "002236", "大成中证360互联网+大数据100指数A",     "1.03", "92.17", "0.98", "0.0101", 9),
    @(4, "010777", "浙商智选家居股票A",                   "0.12", "90.64", "7.78", "0.0093", 2),
    @(5, "519987", "长信恒利优势混合",                     "0.21", "87.52", "4.26", "0.0089", 7),
    @(6, "010778", "浙商智选家居股票C",                   "0.08", "90.64", "7.78", "0.0062", 2),
    @(7, "015774", "长信多利灵活配置混合E",               "0.07", "89.23", "3.80", "0.0027", 10),
    @(8, "013488", "长信多利灵活配置混合C",               "0.05", "89.23", "3.80", "0.0019", 10)
)

$row = 2
foreach ($d in $data) {
    $newSheet.Range("A$row").Value = $d[0]

    # Fund code must stay textual so leading zeros are preserved.
    $bcell = $newSheet.Range("B$row")
    $bcell.NumberFormat = "@"
    $bcell.Value = $d[1]
    $bcell.ClearFormats()

    $newSheet.Range("C$row").Value = $d[2]

    # Numeric-looking columns D:G are stored as plain text in this
    # workbook, not as numbers - force text storage the same way.
    $dcell = $newSheet.Range("D$row")
    $dcell.NumberFormat = "@"
    $dcell.Value = $d[3]
    $dcell.ClearFormats()

    $ecell = $newSheet.Range("E$row")
    $ecell.NumberFormat = "@"
    $ecell.Value = $d[4]
    $ecell.ClearFormats()

    $fcell = $newSheet.Range("F$row")
    $fcell.NumberFormat = "@"
    $fcell.Value = $d[5]
    $fcell.ClearFormats()

    $gcell = $newSheet.Range("G$row")
    $gcell.NumberFormat = "@"
    $gcell.Value = $d[6]
    $gcell.ClearFormats()

    $newSheet.Range("H$row").Value = $d[7]

    $row++
}

Write-Host "2022-Q3 sheet inserted and 总计 sheet updated."
